$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.769.01'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.067.35'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '513.65'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.67'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.434'
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.32'
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.377'
$ws.Range('E11').Value = '  +2.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.593.66'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('E13').Value = '  +1.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.48'
$ws.Range('E14').Value = '  +5.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000163'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.834.57'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.16'
$ws.Range('E17').Value = '  +4.31%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.065.12'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.80'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.06'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '332.10'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.500'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.81'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0901'
$ws.Range('E27').Value = '  -2.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.43'
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.15'
$ws.Range('E29').Value = '  +4.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.80'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.19'
$ws.Range('E31').Value = '  +3.02%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.77'
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '155.20'
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.63'
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.00'
$ws.Range('E36').Value = '  +3.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.27'
$ws.Range('E37').Value = '  +4.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0674'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.109.59'
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.88'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.50'
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.654'
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.273.24'
$ws.Range('E44').Value = '  +2.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0254'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.56'
$ws.Range('E47').Value = '  +4.21%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.935'
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.91'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.727'
$ws.Range('E50').Value = '  +6.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '255.26'
$ws.Range('E51').Value = '  +9.06%  '
